# Insert 9 new trading-day rows (2019-11-18 .. 2019-11-28) right after the
# existing 2019-11-15 row (row 486) and before the existing 2019-11-29 row
# (old row 487, which shifts down to row 496). This grows the sheet from
# A1:I559 to A1:I568.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 487..559 down by 9 rows, leaving 487..495 blank for the new data.
$ws.Range("A487:A495").EntireRow.Insert()

$newRows = @(
    @(487, 1574035200, "2019-11-18", 0.244, 0.244, 0.238, 0.244, 1194500),
    @(488, 1574121600, "2019-11-19", 0.252, 0.274, 0.248, 0.27,  12739250),
    @(489, 1574208000, "2019-11-20", 0.27,  0.27,  0.256, 0.264, 6632750),
    @(490, 1574294400, "2019-11-21", 0.264, 0.264, 0.256, 0.256, 904750),
    @(491, 1574380800, "2019-11-22", 0.26,  0.266, 0.25,  0.25,  2047500),
    @(492, 1574640000, "2019-11-25", 0.248, 0.248, 0.246, 0.246, 703750),
    @(493, 1574726400, "2019-11-26", 0.244, 0.25,  0.244, 0.248, 291750),
    @(494, 1574812800, "2019-11-27", 0.246, 0.25,  0.244, 0.246, 736250),
    @(495, 1574899200, "2019-11-28", 0.246, 0.254, 0.244, 0.246, 587750)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    # Force text so Excel doesn't reinterpret the ISO-looking date string or
    # the zero-padded numeric id as a date/number.
    $ws.Range("B$r").Value = "'" + $row[2]
    $ws.Range("C$r").Value = "'0193"
    $ws.Range("D$r").Value = "KAB"
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $ws.Range("H$r").Value = $row[6]
    $ws.Range("I$r").Value = $row[7]
}

Write-Host "Inserted $($newRows.Count) new rows (487-495)."
